$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 9844.520545567508
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9856.308184707115

$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 114.8270160096505
$ws.Range("D3").Value = 9844.520545567508
$ws.Range("E3").Value = 616238.5361209477
$ws.Range("G3").Value = 626198.5590126801
